# "More colleges input into database"
#
# Appends additional college / program rows to the Colleges sheet and
# reshapes row 3 (UC Berkeley Extension -> UC Berkeley Seminars, program
# not yet filled in). Column C (program name) carries the existing
# "Helvetica 13" cell style (the same style already used by C1/C2/...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# r,  A,    B (school),               C (program),                                                          D (logo/space),                                       E (region)
$rows = @(
    @(1,  1,  "UC Berkeley Extension", "Smart Grid Technologies",                                             "http://extension.berkeley.edu/images/home/logo.jpg", "Northern"),
    @(2,  2,  "UC Berkeley Extension", $null,                                                                 "http://extension.berkeley.edu/images/home/logo.jpg", "Northern"),
    @(3,  3,  "UC Berkeley Seminars",  $null,                                                                 " ",                                                   "Northern"),
    @(4,  4,  "UC Berkeley",           "ME",                                                                  " ",                                                   "Northern"),
    @(5,  5,  "UC Berkeley",           "EE",                                                                  " ",                                                   "Northern"),
    @(6,  6,  "UC Berkeley",           "CompSci",                                                             " ",                                                   "Northern"),
    @(7,  7,  "UC Berkeley",           "Industrial Engineering",                                              " ",                                                   "Northern"),
    @(8,  8,  "UC Berkeley",           "Energy and Resources Group Courses",                                  " ",                                                   "Northern"),
    @(9,  9,  "UC Davis",              "Certification Program in Green Building and Sustainable Design",      " ",                                                   "Northern"),
    @(10, 10, "UC Davis",              "Certificate Program in Energy Resource Management",                   " ",                                                   $null),
    @(11, 11, "UC Davis",              "Certificate Program in Renewable Energy",                             " ",                                                   $null),
    @(12, 12, "UC Davis",              "Professional Concentration in Solar Energy Systems and Design",        " ",                                                   $null),
    @(13, 13, "UC Davis",              "Mechanical Engineering",                                              " ",                                                   $null),
    @(14, 14, "UC Davis",              "Electrical Engineering",                                              " ",                                                   $null),
    @(15, 15, "UC Davis",              "Computer Science",                                                    " ",                                                   $null),
    @(16, 16, "UC Davis",              "Environment Resources Sciences",                                      " ",                                                   $null),
    @(17, 17, "UC Irvine",             "Sustainable Business Management",                                     " ",                                                   $null),
    @(18, 18, "UC Irvine",             "Mechanical Engineering",                                              " ",                                                   $null),
    @(19, 19, "UC Irvine",             "Electrical Engineering",                                              " ",                                                   $null),
    @(20, 20, "UC Irvine",             "Computer Science",                                                    " ",                                                   $null),
    @(21, 21, "UC Irvine",             "Others",                                                               " ",                                                   $null),
    @(22, $null, $null,               $null,                                                                   " ",                                                   $null),
    @(23, $null, $null,               $null,                                                                   " ",                                                   $null),
    @(24, $null, $null,               $null,                                                                   " ",                                                   $null),
    @(25, $null, $null,               $null,                                                                   " ",                                                   $null),
    @(26, $null, $null,               $null,                                                                   " ",                                                   $null),
    @(27, $null, $null,               $null,                                                                   " ",                                                   $null),
    @(28, $null, $null,               $null,                                                                   " ",                                                   $null),
    @(29, $null, $null,               $null,                                                                   " ",                                                   $null),
    @(30, $null, $null,               $null,                                                                   " ",                                                   $null),
    @(31, $null, $null,               $null,                                                                   " ",                                                   $null),
    @(32, $null, $null,               $null,                                                                   " ",                                                   $null),
    @(33, $null, $null,               $null,                                                                   " ",                                                   $null)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $aVal = $r[1]
    $bVal = $r[2]
    $cVal = $r[3]
    $dVal = $r[4]
    $eVal = $r[5]

    if ($aVal -ne $null) {
        $ws.Cells.Item($rowNum, 1).Value = $aVal
    } else {
        $ws.Cells.Item($rowNum, 1).Clear()
    }

    if ($bVal -ne $null) {
        $ws.Cells.Item($rowNum, 2).Value = $bVal
    } else {
        $ws.Cells.Item($rowNum, 2).Clear()
    }

    $cCell = $ws.Cells.Item($rowNum, 3)
    if ($cVal -ne $null) {
        $cCell.Value = $cVal
        $cCell.Font.Name = "Helvetica"
        $cCell.Font.Size = 13
    } else {
        $cCell.Clear()
    }

    if ($dVal -ne $null) {
        $ws.Cells.Item($rowNum, 4).Value = $dVal
    } else {
        $ws.Cells.Item($rowNum, 4).Clear()
    }

    if ($eVal -ne $null) {
        $ws.Cells.Item($rowNum, 5).Value = $eVal
    } else {
        $ws.Cells.Item($rowNum, 5).Clear()
    }
}

# Row 2's program cell (C2) stays blank but keeps the Helvetica style, same
# as before the edit.
$ws.Cells.Item(2, 3).Font.Name = "Helvetica"
$ws.Cells.Item(2, 3).Font.Size = 13

# Row 3 loses its program entirely (no value, no style) - UC Berkeley
# Seminars hasn't had its course filled in yet.
$ws.Cells.Item(3, 3).Clear()

$ws.Range("D16").Select() | Out-Null
